# Clarify cascading select example (#207)
#
# Renames the generic "cf" choice_filter pattern to explicit per-level
# selected_* / state / county columns, trims the duplicate
# dumont/finney choices, and gives the form a friendlier title + id.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("B2").Value = "selected_state"
$survey.Range("C2").Value = "Select a state"

$survey.Range("B3").Value = "selected_county"
$survey.Range("C3").Value = "Select a county"
$survey.Range("D3").Value = "state=`${selected_state}"

$survey.Range("B4").Value = "selected_city"
$survey.Range("C4").Value = "Select a city"
$survey.Range("D4").Value = "county=`${selected_county}"

# ---------------------------------------------------------------
# choices sheet
# ---------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# header: cf -> state, plus new county column
$choices.Range("D1").Value = "state"
$choices.Range("D1").Copy() | Out-Null
$choices.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$choices.Range("E1").Value = "county"

# counties block moves up one row (was rows 4-7, now rows 5-8);
# values themselves are unchanged (king/pierce/brewster/cameron)
$choices.Range("A5").Value = "counties"
$choices.Range("B5").Value = "king"
$choices.Range("C5").Value = "King"
$choices.Range("D5").Value = "washington"

$choices.Range("A6").Value = "counties"
$choices.Range("B6").Value = "pierce"
$choices.Range("C6").Value = "Pierce"
$choices.Range("D6").Value = "washington"

$choices.Range("A7").Value = "counties"
$choices.Range("B7").Value = "brewster"
$choices.Range("C7").Value = "Brewster"
$choices.Range("D7").Value = "texas"

$choices.Range("A8").Value = "counties"
$choices.Range("B8").Value = "cameron"
$choices.Range("C8").Value = "Cameron"
$choices.Range("D8").Value = "texas"

# old row 4 (king) now sits blank below the states block
$choices.Range("A4:D4").Value = ""

# dumont / finney cities are dropped entirely - row 9 is now blank
$choices.Range("A9:D9").Value = ""

# remaining cities (rows 10-15) keep their row position; name/label
# stay put, but the old single "cf" value is replaced by explicit
# state + county columns. The E column is new for these rows, so its
# style (s="6") is pulled from the existing D cell on the same row
# before the value is written.
$choices.Range("A10").Value = "cities"
$choices.Range("D10").Value = "texas"
$choices.Range("D10").Copy() | Out-Null
$choices.Range("E10").PasteSpecial(-4122) | Out-Null
$choices.Range("E10").Value = "cameron"

$choices.Range("A11").Value = "cities"
$choices.Range("D11").Value = "texas"
$choices.Range("D11").Copy() | Out-Null
$choices.Range("E11").PasteSpecial(-4122) | Out-Null
$choices.Range("E11").Value = "cameron"

$choices.Range("A12").Value = "cities"
$choices.Range("D12").Value = "washington"
$choices.Range("D12").Copy() | Out-Null
$choices.Range("E12").PasteSpecial(-4122) | Out-Null
$choices.Range("E12").Value = "king"

$choices.Range("A13").Value = "cities"
$choices.Range("D13").Value = "washington"
$choices.Range("D13").Copy() | Out-Null
$choices.Range("E13").PasteSpecial(-4122) | Out-Null
$choices.Range("E13").Value = "king"

$choices.Range("A14").Value = "cities"
$choices.Range("D14").Value = "washington"
$choices.Range("D14").Copy() | Out-Null
$choices.Range("E14").PasteSpecial(-4122) | Out-Null
$choices.Range("E14").Value = "pierce"

$choices.Range("A15").Value = "cities"
$choices.Range("D15").Value = "washington"
$choices.Range("D15").Copy() | Out-Null
$choices.Range("E15").PasteSpecial(-4122) | Out-Null
$choices.Range("E15").Value = "pierce"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# settings sheet
# ---------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

$settings.Range("A2").Value = "Cascading select example"
$settings.Range("B2").Value = "cascading_select"
